$d = $word.ActiveDocument

# Phase 1: replace each original expression with a unique placeholder token.
# Using placeholders (instead of replacing old -> new directly) avoids collisions
# where one expression's text is a substring of another expression's old or new
# text (e.g. "9+5=" is a substring of "39+5=", and "1+76=" is a substring of
# "21+76="). Longer "old" strings are replaced first so a shorter pattern cannot
# accidentally match inside a longer one that has not been processed yet.
$d.Content.Find.Execute("71-52=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH0@@", 2) | Out-Null
$d.Content.Find.Execute("87-62=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH1@@", 2) | Out-Null
$d.Content.Find.Execute("25+20=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH2@@", 2) | Out-Null
$d.Content.Find.Execute("86-38=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH4@@", 2) | Out-Null
$d.Content.Find.Execute("42-21=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH6@@", 2) | Out-Null
$d.Content.Find.Execute("12+68=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH7@@", 2) | Out-Null
$d.Content.Find.Execute("36+51=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH10@@", 2) | Out-Null
$d.Content.Find.Execute("40+51=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH12@@", 2) | Out-Null
$d.Content.Find.Execute("83-47=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH13@@", 2) | Out-Null
$d.Content.Find.Execute("28+25=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH14@@", 2) | Out-Null
$d.Content.Find.Execute("71+21=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH15@@", 2) | Out-Null
$d.Content.Find.Execute("83+16=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH16@@", 2) | Out-Null
$d.Content.Find.Execute("85-74=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH18@@", 2) | Out-Null
$d.Content.Find.Execute("54-35=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH19@@", 2) | Out-Null
$d.Content.Find.Execute("92-60=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH20@@", 2) | Out-Null
$d.Content.Find.Execute("98-24=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH21@@", 2) | Out-Null
$d.Content.Find.Execute("44-15=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH22@@", 2) | Out-Null
$d.Content.Find.Execute("97-16=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH23@@", 2) | Out-Null
$d.Content.Find.Execute("35+63=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH25@@", 2) | Out-Null
$d.Content.Find.Execute("45+50=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH26@@", 2) | Out-Null
$d.Content.Find.Execute("49+44=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH29@@", 2) | Out-Null
$d.Content.Find.Execute("19+77=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH30@@", 2) | Out-Null
$d.Content.Find.Execute("45-22=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH31@@", 2) | Out-Null
$d.Content.Find.Execute("59-32=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH33@@", 2) | Out-Null
$d.Content.Find.Execute("17+30=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH34@@", 2) | Out-Null
$d.Content.Find.Execute("73-23=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH38@@", 2) | Out-Null
$d.Content.Find.Execute("24+43=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH39@@", 2) | Out-Null
$d.Content.Find.Execute("86-10=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH40@@", 2) | Out-Null
$d.Content.Find.Execute("95-46=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH41@@", 2) | Out-Null
$d.Content.Find.Execute("37+15=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH43@@", 2) | Out-Null
$d.Content.Find.Execute("17+10=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH44@@", 2) | Out-Null
$d.Content.Find.Execute("95-14=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH45@@", 2) | Out-Null
$d.Content.Find.Execute("46+35=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH47@@", 2) | Out-Null
$d.Content.Find.Execute("44+29=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH48@@", 2) | Out-Null
$d.Content.Find.Execute("56-30=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH50@@", 2) | Out-Null
$d.Content.Find.Execute("97-81=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH51@@", 2) | Out-Null
$d.Content.Find.Execute("23-16=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH53@@", 2) | Out-Null
$d.Content.Find.Execute("14+44=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH55@@", 2) | Out-Null
$d.Content.Find.Execute("32-16=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH57@@", 2) | Out-Null
$d.Content.Find.Execute("45+37=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH58@@", 2) | Out-Null
$d.Content.Find.Execute("13+13=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH59@@", 2) | Out-Null
$d.Content.Find.Execute("96-78=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH60@@", 2) | Out-Null
$d.Content.Find.Execute("90-38=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH61@@", 2) | Out-Null
$d.Content.Find.Execute("83+10=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH62@@", 2) | Out-Null
$d.Content.Find.Execute("68+23=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH63@@", 2) | Out-Null
$d.Content.Find.Execute("77-61=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH65@@", 2) | Out-Null
$d.Content.Find.Execute("21+76=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH66@@", 2) | Out-Null
$d.Content.Find.Execute("99-57=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH67@@", 2) | Out-Null
$d.Content.Find.Execute("89-18=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH68@@", 2) | Out-Null
$d.Content.Find.Execute("24-11=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH69@@", 2) | Out-Null
$d.Content.Find.Execute("98-29=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH70@@", 2) | Out-Null
$d.Content.Find.Execute("17+46=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH71@@", 2) | Out-Null
$d.Content.Find.Execute("49+28=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH74@@", 2) | Out-Null
$d.Content.Find.Execute("39+32=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH76@@", 2) | Out-Null
$d.Content.Find.Execute("82-56=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH77@@", 2) | Out-Null
$d.Content.Find.Execute("14+60=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH80@@", 2) | Out-Null
$d.Content.Find.Execute("25+49=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH81@@", 2) | Out-Null
$d.Content.Find.Execute("24+26=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH82@@", 2) | Out-Null
$d.Content.Find.Execute("85-57=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH83@@", 2) | Out-Null
$d.Content.Find.Execute("27-19=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH84@@", 2) | Out-Null
$d.Content.Find.Execute("28+49=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH85@@", 2) | Out-Null
$d.Content.Find.Execute("63+27=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH86@@", 2) | Out-Null
$d.Content.Find.Execute("62+34=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH88@@", 2) | Out-Null
$d.Content.Find.Execute("11+62=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH90@@", 2) | Out-Null
$d.Content.Find.Execute("41-34=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH91@@", 2) | Out-Null
$d.Content.Find.Execute("68+24=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH93@@", 2) | Out-Null
$d.Content.Find.Execute("98-62=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH94@@", 2) | Out-Null
$d.Content.Find.Execute("88-83=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH95@@", 2) | Out-Null
$d.Content.Find.Execute("11+57=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH96@@", 2) | Out-Null
$d.Content.Find.Execute("41+40=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH97@@", 2) | Out-Null
$d.Content.Find.Execute("20+67=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH98@@", 2) | Out-Null
$d.Content.Find.Execute("2+69=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH3@@", 2) | Out-Null
$d.Content.Find.Execute("82+3=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH5@@", 2) | Out-Null
$d.Content.Find.Execute("1+76=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH8@@", 2) | Out-Null
$d.Content.Find.Execute("35+0=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH9@@", 2) | Out-Null
$d.Content.Find.Execute("6+88=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH11@@", 2) | Out-Null
$d.Content.Find.Execute("73+7=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH17@@", 2) | Out-Null
$d.Content.Find.Execute("70+3=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH24@@", 2) | Out-Null
$d.Content.Find.Execute("59+1=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH27@@", 2) | Out-Null
$d.Content.Find.Execute("6+37=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH28@@", 2) | Out-Null
$d.Content.Find.Execute("4+68=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH32@@", 2) | Out-Null
$d.Content.Find.Execute("23+0=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH35@@", 2) | Out-Null
$d.Content.Find.Execute("21+9=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH36@@", 2) | Out-Null
$d.Content.Find.Execute("8+88=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH37@@", 2) | Out-Null
$d.Content.Find.Execute("81-3=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH46@@", 2) | Out-Null
$d.Content.Find.Execute("2+75=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH49@@", 2) | Out-Null
$d.Content.Find.Execute("13+7=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH52@@", 2) | Out-Null
$d.Content.Find.Execute("45-9=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH54@@", 2) | Out-Null
$d.Content.Find.Execute("8+41=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH56@@", 2) | Out-Null
$d.Content.Find.Execute("47+3=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH64@@", 2) | Out-Null
$d.Content.Find.Execute("52+2=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH72@@", 2) | Out-Null
$d.Content.Find.Execute("54-8=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH73@@", 2) | Out-Null
$d.Content.Find.Execute("93+1=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH75@@", 2) | Out-Null
$d.Content.Find.Execute("8+19=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH78@@", 2) | Out-Null
$d.Content.Find.Execute("83-4=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH79@@", 2) | Out-Null
$d.Content.Find.Execute("1+71=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH87@@", 2) | Out-Null
$d.Content.Find.Execute("75-6=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH89@@", 2) | Out-Null
$d.Content.Find.Execute("96+3=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH92@@", 2) | Out-Null
$d.Content.Find.Execute("7+4=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH42@@", 2) | Out-Null
$d.Content.Find.Execute("9+5=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH99@@", 2) | Out-Null

# Phase 2: replace each placeholder with the final new expression.
$d.Content.Find.Execute("@@PH0@@", $true, $false, $false, $false, $false, $true, 1, $false, "47+24=", 2) | Out-Null
$d.Content.Find.Execute("@@PH1@@", $true, $false, $false, $false, $false, $true, 1, $false, "32+5=", 2) | Out-Null
$d.Content.Find.Execute("@@PH2@@", $true, $false, $false, $false, $false, $true, 1, $false, "35+7=", 2) | Out-Null
$d.Content.Find.Execute("@@PH3@@", $true, $false, $false, $false, $false, $true, 1, $false, "58-3=", 2) | Out-Null
$d.Content.Find.Execute("@@PH4@@", $true, $false, $false, $false, $false, $true, 1, $false, "57+2=", 2) | Out-Null
$d.Content.Find.Execute("@@PH5@@", $true, $false, $false, $false, $false, $true, 1, $false, "88-84=", 2) | Out-Null
$d.Content.Find.Execute("@@PH6@@", $true, $false, $false, $false, $false, $true, 1, $false, "79-60=", 2) | Out-Null
$d.Content.Find.Execute("@@PH7@@", $true, $false, $false, $false, $false, $true, 1, $false, "22+57=", 2) | Out-Null
$d.Content.Find.Execute("@@PH8@@", $true, $false, $false, $false, $false, $true, 1, $false, "87-39=", 2) | Out-Null
$d.Content.Find.Execute("@@PH9@@", $true, $false, $false, $false, $false, $true, 1, $false, "9+58=", 2) | Out-Null
$d.Content.Find.Execute("@@PH10@@", $true, $false, $false, $false, $false, $true, 1, $false, "39-33=", 2) | Out-Null
$d.Content.Find.Execute("@@PH11@@", $true, $false, $false, $false, $false, $true, 1, $false, "42-4=", 2) | Out-Null
$d.Content.Find.Execute("@@PH12@@", $true, $false, $false, $false, $false, $true, 1, $false, "11+46=", 2) | Out-Null
$d.Content.Find.Execute("@@PH13@@", $true, $false, $false, $false, $false, $true, 1, $false, "17+5=", 2) | Out-Null
$d.Content.Find.Execute("@@PH14@@", $true, $false, $false, $false, $false, $true, 1, $false, "52+0=", 2) | Out-Null
$d.Content.Find.Execute("@@PH15@@", $true, $false, $false, $false, $false, $true, 1, $false, "64-14=", 2) | Out-Null
$d.Content.Find.Execute("@@PH16@@", $true, $false, $false, $false, $false, $true, 1, $false, "5+41=", 2) | Out-Null
$d.Content.Find.Execute("@@PH17@@", $true, $false, $false, $false, $false, $true, 1, $false, "39+27=", 2) | Out-Null
$d.Content.Find.Execute("@@PH18@@", $true, $false, $false, $false, $false, $true, 1, $false, "51+23=", 2) | Out-Null
$d.Content.Find.Execute("@@PH19@@", $true, $false, $false, $false, $false, $true, 1, $false, "78-32=", 2) | Out-Null
$d.Content.Find.Execute("@@PH20@@", $true, $false, $false, $false, $false, $true, 1, $false, "57+19=", 2) | Out-Null
$d.Content.Find.Execute("@@PH21@@", $true, $false, $false, $false, $false, $true, 1, $false, "7+73=", 2) | Out-Null
$d.Content.Find.Execute("@@PH22@@", $true, $false, $false, $false, $false, $true, 1, $false, "0+42=", 2) | Out-Null
$d.Content.Find.Execute("@@PH23@@", $true, $false, $false, $false, $false, $true, 1, $false, "92-21=", 2) | Out-Null
$d.Content.Find.Execute("@@PH24@@", $true, $false, $false, $false, $false, $true, 1, $false, "22+3=", 2) | Out-Null
$d.Content.Find.Execute("@@PH25@@", $true, $false, $false, $false, $false, $true, 1, $false, "70-47=", 2) | Out-Null
$d.Content.Find.Execute("@@PH26@@", $true, $false, $false, $false, $false, $true, 1, $false, "29+34=", 2) | Out-Null
$d.Content.Find.Execute("@@PH27@@", $true, $false, $false, $false, $false, $true, 1, $false, "56-1=", 2) | Out-Null
$d.Content.Find.Execute("@@PH28@@", $true, $false, $false, $false, $false, $true, 1, $false, "64-64=", 2) | Out-Null
$d.Content.Find.Execute("@@PH29@@", $true, $false, $false, $false, $false, $true, 1, $false, "72+20=", 2) | Out-Null
$d.Content.Find.Execute("@@PH30@@", $true, $false, $false, $false, $false, $true, 1, $false, "84-44=", 2) | Out-Null
$d.Content.Find.Execute("@@PH31@@", $true, $false, $false, $false, $false, $true, 1, $false, "43-41=", 2) | Out-Null
$d.Content.Find.Execute("@@PH32@@", $true, $false, $false, $false, $false, $true, 1, $false, "91-58=", 2) | Out-Null
$d.Content.Find.Execute("@@PH33@@", $true, $false, $false, $false, $false, $true, 1, $false, "59-26=", 2) | Out-Null
$d.Content.Find.Execute("@@PH34@@", $true, $false, $false, $false, $false, $true, 1, $false, "54-0=", 2) | Out-Null
$d.Content.Find.Execute("@@PH35@@", $true, $false, $false, $false, $false, $true, 1, $false, "87-81=", 2) | Out-Null
$d.Content.Find.Execute("@@PH36@@", $true, $false, $false, $false, $false, $true, 1, $false, "13+25=", 2) | Out-Null
$d.Content.Find.Execute("@@PH37@@", $true, $false, $false, $false, $false, $true, 1, $false, "56+9=", 2) | Out-Null
$d.Content.Find.Execute("@@PH38@@", $true, $false, $false, $false, $false, $true, 1, $false, "17+78=", 2) | Out-Null
$d.Content.Find.Execute("@@PH39@@", $true, $false, $false, $false, $false, $true, 1, $false, "58-48=", 2) | Out-Null
$d.Content.Find.Execute("@@PH40@@", $true, $false, $false, $false, $false, $true, 1, $false, "45+16=", 2) | Out-Null
$d.Content.Find.Execute("@@PH41@@", $true, $false, $false, $false, $false, $true, 1, $false, "41-39=", 2) | Out-Null
$d.Content.Find.Execute("@@PH42@@", $true, $false, $false, $false, $false, $true, 1, $false, "44-42=", 2) | Out-Null
$d.Content.Find.Execute("@@PH43@@", $true, $false, $false, $false, $false, $true, 1, $false, "48+7=", 2) | Out-Null
$d.Content.Find.Execute("@@PH44@@", $true, $false, $false, $false, $false, $true, 1, $false, "43-3=", 2) | Out-Null
$d.Content.Find.Execute("@@PH45@@", $true, $false, $false, $false, $false, $true, 1, $false, "80-63=", 2) | Out-Null
$d.Content.Find.Execute("@@PH46@@", $true, $false, $false, $false, $false, $true, 1, $false, "29+49=", 2) | Out-Null
$d.Content.Find.Execute("@@PH47@@", $true, $false, $false, $false, $false, $true, 1, $false, "21-13=", 2) | Out-Null
$d.Content.Find.Execute("@@PH48@@", $true, $false, $false, $false, $false, $true, 1, $false, "16+7=", 2) | Out-Null
$d.Content.Find.Execute("@@PH49@@", $true, $false, $false, $false, $false, $true, 1, $false, "39+5=", 2) | Out-Null
$d.Content.Find.Execute("@@PH50@@", $true, $false, $false, $false, $false, $true, 1, $false, "14+31=", 2) | Out-Null
$d.Content.Find.Execute("@@PH51@@", $true, $false, $false, $false, $false, $true, 1, $false, "91-38=", 2) | Out-Null
$d.Content.Find.Execute("@@PH52@@", $true, $false, $false, $false, $false, $true, 1, $false, "90-34=", 2) | Out-Null
$d.Content.Find.Execute("@@PH53@@", $true, $false, $false, $false, $false, $true, 1, $false, "16+57=", 2) | Out-Null
$d.Content.Find.Execute("@@PH54@@", $true, $false, $false, $false, $false, $true, 1, $false, "21+1=", 2) | Out-Null
$d.Content.Find.Execute("@@PH55@@", $true, $false, $false, $false, $false, $true, 1, $false, "43+53=", 2) | Out-Null
$d.Content.Find.Execute("@@PH56@@", $true, $false, $false, $false, $false, $true, 1, $false, "20+5=", 2) | Out-Null
$d.Content.Find.Execute("@@PH57@@", $true, $false, $false, $false, $false, $true, 1, $false, "23-8=", 2) | Out-Null
$d.Content.Find.Execute("@@PH58@@", $true, $false, $false, $false, $false, $true, 1, $false, "97-34=", 2) | Out-Null
$d.Content.Find.Execute("@@PH59@@", $true, $false, $false, $false, $false, $true, 1, $false, "12+32=", 2) | Out-Null
$d.Content.Find.Execute("@@PH60@@", $true, $false, $false, $false, $false, $true, 1, $false, "72+21=", 2) | Out-Null
$d.Content.Find.Execute("@@PH61@@", $true, $false, $false, $false, $false, $true, 1, $false, "83-46=", 2) | Out-Null
$d.Content.Find.Execute("@@PH62@@", $true, $false, $false, $false, $false, $true, 1, $false, "6+80=", 2) | Out-Null
$d.Content.Find.Execute("@@PH63@@", $true, $false, $false, $false, $false, $true, 1, $false, "11+1=", 2) | Out-Null
$d.Content.Find.Execute("@@PH64@@", $true, $false, $false, $false, $false, $true, 1, $false, "80-48=", 2) | Out-Null
$d.Content.Find.Execute("@@PH65@@", $true, $false, $false, $false, $false, $true, 1, $false, "45+28=", 2) | Out-Null
$d.Content.Find.Execute("@@PH66@@", $true, $false, $false, $false, $false, $true, 1, $false, "91+4=", 2) | Out-Null
$d.Content.Find.Execute("@@PH67@@", $true, $false, $false, $false, $false, $true, 1, $false, "60-18=", 2) | Out-Null
$d.Content.Find.Execute("@@PH68@@", $true, $false, $false, $false, $false, $true, 1, $false, "30+9=", 2) | Out-Null
$d.Content.Find.Execute("@@PH69@@", $true, $false, $false, $false, $false, $true, 1, $false, "77-23=", 2) | Out-Null
$d.Content.Find.Execute("@@PH70@@", $true, $false, $false, $false, $false, $true, 1, $false, "69-25=", 2) | Out-Null
$d.Content.Find.Execute("@@PH71@@", $true, $false, $false, $false, $false, $true, 1, $false, "23+25=", 2) | Out-Null
$d.Content.Find.Execute("@@PH72@@", $true, $false, $false, $false, $false, $true, 1, $false, "26+66=", 2) | Out-Null
$d.Content.Find.Execute("@@PH73@@", $true, $false, $false, $false, $false, $true, 1, $false, "16+20=", 2) | Out-Null
$d.Content.Find.Execute("@@PH74@@", $true, $false, $false, $false, $false, $true, 1, $false, "65-40=", 2) | Out-Null
$d.Content.Find.Execute("@@PH75@@", $true, $false, $false, $false, $false, $true, 1, $false, "88-35=", 2) | Out-Null
$d.Content.Find.Execute("@@PH76@@", $true, $false, $false, $false, $false, $true, 1, $false, "70+5=", 2) | Out-Null
$d.Content.Find.Execute("@@PH77@@", $true, $false, $false, $false, $false, $true, 1, $false, "62-32=", 2) | Out-Null
$d.Content.Find.Execute("@@PH78@@", $true, $false, $false, $false, $false, $true, 1, $false, "76-37=", 2) | Out-Null
$d.Content.Find.Execute("@@PH79@@", $true, $false, $false, $false, $false, $true, 1, $false, "27+52=", 2) | Out-Null
$d.Content.Find.Execute("@@PH80@@", $true, $false, $false, $false, $false, $true, 1, $false, "16-16=", 2) | Out-Null
$d.Content.Find.Execute("@@PH81@@", $true, $false, $false, $false, $false, $true, 1, $false, "74-1=", 2) | Out-Null
$d.Content.Find.Execute("@@PH82@@", $true, $false, $false, $false, $false, $true, 1, $false, "7+30=", 2) | Out-Null
$d.Content.Find.Execute("@@PH83@@", $true, $false, $false, $false, $false, $true, 1, $false, "20+31=", 2) | Out-Null
$d.Content.Find.Execute("@@PH84@@", $true, $false, $false, $false, $false, $true, 1, $false, "97-82=", 2) | Out-Null
$d.Content.Find.Execute("@@PH85@@", $true, $false, $false, $false, $false, $true, 1, $false, "14+11=", 2) | Out-Null
$d.Content.Find.Execute("@@PH86@@", $true, $false, $false, $false, $false, $true, 1, $false, "94-13=", 2) | Out-Null
$d.Content.Find.Execute("@@PH87@@", $true, $false, $false, $false, $false, $true, 1, $false, "10+60=", 2) | Out-Null
$d.Content.Find.Execute("@@PH88@@", $true, $false, $false, $false, $false, $true, 1, $false, "77+22=", 2) | Out-Null
$d.Content.Find.Execute("@@PH89@@", $true, $false, $false, $false, $false, $true, 1, $false, "40-29=", 2) | Out-Null
$d.Content.Find.Execute("@@PH90@@", $true, $false, $false, $false, $false, $true, 1, $false, "40+36=", 2) | Out-Null
$d.Content.Find.Execute("@@PH91@@", $true, $false, $false, $false, $false, $true, 1, $false, "90+2=", 2) | Out-Null
$d.Content.Find.Execute("@@PH92@@", $true, $false, $false, $false, $false, $true, 1, $false, "7+3=", 2) | Out-Null
$d.Content.Find.Execute("@@PH93@@", $true, $false, $false, $false, $false, $true, 1, $false, "94-15=", 2) | Out-Null
$d.Content.Find.Execute("@@PH94@@", $true, $false, $false, $false, $false, $true, 1, $false, "29+44=", 2) | Out-Null
$d.Content.Find.Execute("@@PH95@@", $true, $false, $false, $false, $false, $true, 1, $false, "51+14=", 2) | Out-Null
$d.Content.Find.Execute("@@PH96@@", $true, $false, $false, $false, $false, $true, 1, $false, "65-30=", 2) | Out-Null
$d.Content.Find.Execute("@@PH97@@", $true, $false, $false, $false, $false, $true, 1, $false, "9+80=", 2) | Out-Null
$d.Content.Find.Execute("@@PH98@@", $true, $false, $false, $false, $false, $true, 1, $false, "74+14=", 2) | Out-Null
$d.Content.Find.Execute("@@PH99@@", $true, $false, $false, $false, $false, $true, 1, $false, "33+15=", 2) | Out-Null
